$d = $word.ActiveDocument

function Merge-ParagraphRuns($paragraphIndex, $newText) {
    $p = $d.Paragraphs.Item($paragraphIndex).Range
    $p.MoveEnd(1, -1) | Out-Null
    $p.Delete() | Out-Null
    $p.InsertAfter($newText) | Out-Null
}

Merge-ParagraphRuns 1 "Answers: Rearranging equations involving trigonometry and logarithms"
Merge-ParagraphRuns 2 "Ellie Gurini"
Merge-ParagraphRuns 4 "This is an answer set relating to the questions based on Guide, Introduction to rearranging equations involving trigonometry and logarithms."
